# Commit: "added effect of frost and heat on LAI decrease"
#
# This adds two new computed variables (cFrost, cHeat) describing the
# decrease of LAI due to frost/heat to the "savedEachDay" sheet, adds
# explanatory remarks to three existing frost/heat-threshold parameters,
# and removes a stray translationSSM value that had leaked into the
# sDecreaseLAIperBD row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "savedEachDay" is already the active sheet

# --- Remarks (column J) explaining the in-code names of three existing
#     frost/heat parameters -------------------------------------------------
$ws.Range("J64").Value = "in the code this is called FRZLDR"
$ws.Range("J65").Value = "in the code this is called FRZTKIL"
$ws.Range("J67").Value = "in the code this is called HtLTH"

# --- Drop the stray translationSSM value erroneously left on row 116 ------
$ws.Range("G116").ClearContents()

# --- New rows describing the two new computed variables -------------------
# Filled column-by-column (both rows at once) to match the order variables
# were actually typed in.
$ws.Range("A117").Value = "cFrost"
$ws.Range("A118").Value = "cHeat"

$ws.Range("B117").Value = "computed"
$ws.Range("B118").Value = "computed"

$ws.Range("C117").Value = "numeric"
$ws.Range("C118").Value = "numeric"

$ws.Range("D117").Value = "LAI_Senescence"
$ws.Range("D118").Value = "LAI_Senescence"

$ws.Range("E117").Value = "m2 m-2"
$ws.Range("E118").Value = "m2 m-2"

# Copy the "definition" column number format/font from the row above (row
# 116) so the new cells pick up the same style index, then fill the text.
$ws.Range("F116").Copy()
$ws.Range("F117:F118").PasteSpecial(-4122)
$ws.Range("F117").Value = "decrease of LAI due to frost"
$ws.Range("F118").Value = "decrease of LAI due to heat"

$ws.Range("G117").Value = "DLAIF"
$ws.Range("G118").Value = "DLAIH"

$ws.Range("H117").Value = "décroissance du LAI à cause du froid"
$ws.Range("H118").Value = "décroissance du LAI à cause de la chaleur"

$ws.Range("I117").Formula = "=NA()"
$ws.Range("I118").Formula = "=NA()"

# --- Update the on-screen selection to match where the editor ended up ----
[void]$ws.Range("J113").Select()
